$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 12,20
$data[0,0] = "ECs"
$data[0,1] = "Nppa"
$data[0,2] = "Npr3"
$data[0,3] = "ECs"
$data[0,4] = 1
$data[0,5] = 0.3333333333333333
$data[0,6] = 0.110293
$data[0,7] = 0.330879
$data[0,8] = 0.1397687944778318
$data[0,9] = 0.1397687944778318
$data[0,10] = 1
$data[0,11] = 0.3333333333333333
$data[0,12] = 0.1036536666666667
$data[0,13] = 0.310961
$data[0,14] = 0.08145093039891602
$data[0,15] = 0.08145093039891602
$data[0,16] = 0.01143227385766666
$data[0,17] = 0.102890464719
$data[0,18] = 0.01138429835095427
$data[0,19] = 0.01138429835095427
$data[1,0] = "ECs"
$data[1,1] = "Nppa"
$data[1,2] = "Npr3"
$data[1,3] = "FAPs"
$data[1,4] = 1
$data[1,5] = 0.3333333333333333
$data[1,6] = 0.110293
$data[1,7] = 0.330879
$data[1,8] = 0.1397687944778318
$data[1,9] = 0.1397687944778318
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 1.007656
$data[1,13] = 3.022968
$data[1,14] = 0.7918149097994615
$data[1,15] = 0.7918149097994616
$data[1,16] = 0.111137403208
$data[1,17] = 1.000236628872
$data[1,18] = 0.1106710153922438
$data[1,19] = 0.1106710153922438
$data[2,0] = "ECs"
$data[2,1] = "Nppa"
$data[2,2] = "Npr3"
$data[2,3] = "MuSCs"
$data[2,4] = 1
$data[2,5] = 0.3333333333333333
$data[2,6] = 0.110293
$data[2,7] = 0.330879
$data[2,8] = 0.1397687944778318
$data[2,9] = 0.1397687944778318
$data[2,10] = 2
$data[2,11] = 0.6666666666666666
$data[2,12] = 0.1612806666666667
$data[2,13] = 0.483842
$data[2,14] = 0.1267341598016225
$data[2,15] = 0.1267341598016225
$data[2,16] = 0.01778812856866666
$data[2,17] = 0.160093157118
$data[2,18] = 0.01771348073463366
$data[2,19] = 0.01771348073463366
$data[3,0] = "FAPs"
$data[3,1] = "Nppa"
$data[3,2] = "Npr3"
$data[3,3] = "ECs"
$data[3,4] = 1
$data[3,5] = 0.3333333333333333
$data[3,6] = 0.03729466666666666
$data[3,7] = 0.111884
$data[3,8] = 0.04726166302895539
$data[3,9] = 0.04726166302895539
$data[3,10] = 1
$data[3,11] = 0.3333333333333333
$data[3,12] = 0.1036536666666667
$data[3,13] = 0.310961
$data[3,14] = 0.08145093039891602
$data[3,15] = 0.08145093039891602
$data[3,16] = 0.003865728947111111
$data[3,17] = 0.034791560524
$data[3,18] = 0.003849506425908468
$data[3,19] = 0.003849506425908468
$data[4,0] = "FAPs"
$data[4,1] = "Nppa"
$data[4,2] = "Npr3"
$data[4,3] = "FAPs"
$data[4,4] = 1
$data[4,5] = 0.3333333333333333
$data[4,6] = 0.03729466666666666
$data[4,7] = 0.111884
$data[4,8] = 0.04726166302895539
$data[4,9] = 0.04726166302895539
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 1.007656
$data[4,13] = 3.022968
$data[4,14] = 0.7918149097994615
$data[4,15] = 0.7918149097994616
$data[4,16] = 0.03758019463466666
$data[4,17] = 0.338221751712
$data[4,18] = 0.03742248944824485
$data[4,19] = 0.03742248944824486
$data[5,0] = "FAPs"
$data[5,1] = "Nppa"
$data[5,2] = "Npr3"
$data[5,3] = "MuSCs"
$data[5,4] = 1
$data[5,5] = 0.3333333333333333
$data[5,6] = 0.03729466666666666
$data[5,7] = 0.111884
$data[5,8] = 0.04726166302895539
$data[5,9] = 0.04726166302895539
$data[5,10] = 2
$data[5,11] = 0.6666666666666666
$data[5,12] = 0.1612806666666667
$data[5,13] = 0.483842
$data[5,14] = 0.1267341598016225
$data[5,15] = 0.1267341598016225
$data[5,16] = 0.00601490870311111
$data[5,17] = 0.054134178328
$data[5,18] = 0.005989667154802065
$data[5,19] = 0.005989667154802065
$data[6,0] = "MuSCs"
$data[6,1] = "Nppa"
$data[6,2] = "Npr3"
$data[6,3] = "ECs"
$data[6,4] = 2
$data[6,5] = 0.6666666666666666
$data[6,6] = 0.3940343333333334
$data[6,7] = 1.182103
$data[6,8] = 0.4993399740044802
$data[6,9] = 0.4993399740044802
$data[6,10] = 1
$data[6,11] = 0.3333333333333333
$data[6,12] = 0.1036536666666667
$data[6,13] = 0.310961
$data[6,14] = 0.08145093039891602
$data[6,15] = 0.08145093039891602
$data[6,16] = 0.04084310344255555
$data[6,17] = 0.367587930983
$data[6,18] = 0.04067170546803545
$data[6,19] = 0.04067170546803545
$data[7,0] = "MuSCs"
$data[7,1] = "Nppa"
$data[7,2] = "Npr3"
$data[7,3] = "FAPs"
$data[7,4] = 2
$data[7,5] = 0.6666666666666666
$data[7,6] = 0.3940343333333334
$data[7,7] = 1.182103
$data[7,8] = 0.4993399740044802
$data[7,9] = 0.4993399740044802
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 1.007656
$data[7,13] = 3.022968
$data[7,14] = 0.7918149097994615
$data[7,15] = 0.7918149097994616
$data[7,16] = 0.3970510601893333
$data[7,17] = 3.573459541704
$data[7,18] = 0.3953848364756229
$data[7,19] = 0.395384836475623
$data[8,0] = "MuSCs"
$data[8,1] = "Nppa"
$data[8,2] = "Npr3"
$data[8,3] = "MuSCs"
$data[8,4] = 2
$data[8,5] = 0.6666666666666666
$data[8,6] = 0.3940343333333334
$data[8,7] = 1.182103
$data[8,8] = 0.4993399740044802
$data[8,9] = 0.4993399740044802
$data[8,10] = 2
$data[8,11] = 0.6666666666666666
$data[8,12] = 0.1612806666666667
$data[8,13] = 0.483842
$data[8,14] = 0.1267341598016225
$data[8,15] = 0.1267341598016225
$data[8,16] = 0.06355011996955556
$data[8,17] = 0.571951079726
$data[8,18] = 0.0632834320608218
$data[8,19] = 0.0632834320608218
$data[9,0] = "Resolving-Mac"
$data[9,1] = "Nppa"
$data[9,2] = "Npr3"
$data[9,3] = "ECs"
$data[9,4] = 1
$data[9,5] = 0.3333333333333333
$data[9,6] = 0.2474883333333333
$data[9,7] = 0.742465
$data[9,8] = 0.3136295684887327
$data[9,9] = 0.3136295684887327
$data[9,10] = 1
$data[9,11] = 0.3333333333333333
$data[9,12] = 0.1036536666666667
$data[9,13] = 0.310961
$data[9,14] = 0.08145093039891602
$data[9,15] = 0.08145093039891602
$data[9,16] = 0.02565307320722222
$data[9,17] = 0.230877658865
$data[9,18] = 0.02554542015401783
$data[9,19] = 0.02554542015401783
$data[10,0] = "Resolving-Mac"
$data[10,1] = "Nppa"
$data[10,2] = "Npr3"
$data[10,3] = "FAPs"
$data[10,4] = 1
$data[10,5] = 0.3333333333333333
$data[10,6] = 0.2474883333333333
$data[10,7] = 0.742465
$data[10,8] = 0.3136295684887327
$data[10,9] = 0.3136295684887327
$data[10,10] = 3
$data[10,11] = 1
$data[10,12] = 1.007656
$data[10,13] = 3.022968
$data[10,14] = 0.7918149097994615
$data[10,15] = 0.7918149097994616
$data[10,16] = 0.2493831040133333
$data[10,17] = 2.24444793612
$data[10,18] = 0.2483365684833499
$data[10,19] = 0.2483365684833499
$data[11,0] = "Resolving-Mac"
$data[11,1] = "Nppa"
$data[11,2] = "Npr3"
$data[11,3] = "MuSCs"
$data[11,4] = 1
$data[11,5] = 0.3333333333333333
$data[11,6] = 0.2474883333333333
$data[11,7] = 0.742465
$data[11,8] = 0.3136295684887327
$data[11,9] = 0.3136295684887327
$data[11,10] = 2
$data[11,11] = 0.6666666666666666
$data[11,12] = 0.1612806666666667
$data[11,13] = 0.483842
$data[11,14] = 0.1267341598016225
$data[11,15] = 0.1267341598016225
$data[11,16] = 0.03991508339222222
$data[11,17] = 0.35923575053
$data[11,18] = 0.03974757985136495
$data[11,19] = 0.03974757985136495

$ws.Range("A2:T13").Value = $data

$ws.Range("A1").Select() | Out-Null
